$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("truthStateParams")
$ws.Range("B5").NumberFormat = "0.00E+00"
$ws.Range("B5").Font.Bold = $False
$ws.Range("B5").Value = 0.00006
Write-Host "done"
